$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

# ---- Sheet "data": add column Y (25) with header "16. 2. 2021" ----
$wsData.Range("X1").Copy()
$wsData.Range("Y1").PasteSpecial(-4122)
$wsData.Cells.Item(1, 25).Value = "16. 2. 2021"

$yValues = 0.88,0.7,0.65,0.65,0.61,0.34,0.84,0.85,0.93,0.93,0.89,0.87,0.74,0.86,0.95,0.88,0.91,0.85,0.86,0.9,0.88,0.86,0.89,0.9,0.59,0.67,0.79,0.79,0.7,0.65,0.51,0.65,0.81,0.72,0.72,0.65,0.66,0.73,0.7,0.7,0.65,0.7,0.48,0.64,0.76,0.76,0.66,0.59,0.46,0.6,0.77,0.65,0.7,0.61,0.64,0.66,0.63,0.68,0.66,0.67,0.51,0.63,0.74,0.74,0.69,0.58,0.48,0.62,0.74,0.65,0.69,0.6,0.64,0.66,0.63,0.68,0.64,0.64,0.46,0.58,0.72,0.68,0.64,0.58,0.47,0.57,0.7,0.57,0.68,0.58,0.6,0.62,0.57,0.63,0.66,0.68,0.26,0.25,0.48,0.53,0.33,0.21,0.21,0.3,0.43,0.3,0.35,0.37,0.35,0.34,0.37,0.33,0.31,0.3
for ($i = 0; $i -lt $yValues.Length; $i++) {
    $wsData.Cells.Item($i + 2, 25).Value = $yValues[$i]
}

# Update the "aktualizace" date text in the footer row (row 116)
$wsData.Cells.Item(116, 1).Value = "Život během pandemie, Jednotlivé protektivní aktivity, % respondentů celkově a ve skupinách, aktualizace 23. 2. 2021"

# ---- Sheet "pocetR": add column X (24) with header "16. 2. 2021" ----
$wsPocetR.Range("W1").Copy()
$wsPocetR.Range("X1").PasteSpecial(-4122)
$wsPocetR.Cells.Item(1, 24).Value = "16. 2. 2021"

$xValues = 2120,511,774,835,581,703,605,398,744,978,694,668,758,1025,1095,1102,490,248,280
for ($i = 0; $i -lt $xValues.Length; $i++) {
    $wsPocetR.Cells.Item($i + 2, 24).Value = $xValues[$i]
}

# Row 21: trailing empty-string cell under the new column + updated footer text
$wsPocetR.Cells.Item(21, 24).Value = ""
$wsPocetR.Cells.Item(21, 1).Value = "Život během pandemie, Jednotlivé protektivní aktivity, velikost dotázaného souboru celkově a ve skupinách, aktualizace 23. 2. 2021"
